$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Baseline Schedule")

# Update Successors for Activity 1 (row 3): predecessor object repr -> real id 2
$ws.Range("E3").Value = "[(2, 'FS', 0)]"

# Update Predecessors for Activity 2 (row 4): predecessor object repr -> real id 1
$ws.Range("D4").Value = "[(1, 'FS', 0)]"

# Update Baseline Start / Baseline End numeric (serial date) values
$ws.Range("F3").Value = 42103.5082590986
$ws.Range("H3").Value = 42108.5082590988
$ws.Range("F4").Value = 42103.5082590989
$ws.Range("H4").Value = 42113.508259099
